# Adafruit IO feed export: append the newest temperature reading as a new
# row at the bottom of the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 93

# Column C ("Value") holds a numeric-looking reading ("25") but the sheet
# stores every column as text, so force text formatting before assigning it
# to avoid Excel re-typing it as a number.
$ws.Range("C" + $newRow).NumberFormat = "@"

$ws.Range("A" + $newRow).Value = "2024-09-25T18:06:40Z"
$ws.Range("B" + $newRow).Value = "temperature"
$ws.Range("C" + $newRow).Value = "25"
$ws.Range("D" + $newRow).Value = "N/A"
$ws.Range("E" + $newRow).Value = "N/A"
$ws.Range("F" + $newRow).Value = "N/A"
